$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9005614640222286
$ws.Range("J2").Value = 0.9005614640222285
$ws.Range("M2").Value = 36.81180933333333
$ws.Range("N2").Value = 110.435428
$ws.Range("O2").Value = 0.2598784967371026
$ws.Range("P2").Value = 0.2598784967371026
$ws.Range("Q2").Value = 2.883334048445777
$ws.Range("R2").Value = 25.950006436012
$ws.Range("S2").Value = 0.2340365594894611
$ws.Range("T2").Value = 0.2340365594894611
$ws.Range("I3").Value = 0.9005614640222286
$ws.Range("J3").Value = 0.9005614640222285
$ws.Range("O3").Value = 0.1970278712683331
$ws.Range("P3").Value = 0.197027871268333
$ws.Range("S3").Value = 0.1774357082025932
$ws.Range("T3").Value = 0.1774357082025932
$ws.Range("I4").Value = 0.9005614640222286
$ws.Range("J4").Value = 0.9005614640222285
$ws.Range("M4").Value = 21.95609833333333
$ws.Range("N4").Value = 65.868295
$ws.Range("O4").Value = 0.1550023737603119
$ws.Range("P4").Value = 0.1550023737603119
$ws.Range("Q4").Value = 1.719740676756111
$ws.Range("R4").Value = 15.477666090805
$ws.Range("S4").Value = 0.1395891646405072
$ws.Range("T4").Value = 0.1395891646405072
$ws.Range("I5").Value = 0.9005614640222286
$ws.Range("J5").Value = 0.9005614640222285
$ws.Range("M5").Value = 13.23098133333333
$ws.Range("N5").Value = 39.692944
$ws.Range("O5").Value = 0.09340609987756826
$ws.Range("P5").Value = 0.09340609987756825
$ws.Range("Q5").Value = 1.036334254241778
$ws.Range("R5").Value = 9.327008288176
$ws.Range("S5").Value = 0.08411793405434938
$ws.Range("T5").Value = 0.08411793405434935
$ws.Range("I6").Value = 0.9005614640222286
$ws.Range("J6").Value = 0.9005614640222285
$ws.Range("M6").Value = 22.080681
$ws.Range("N6").Value = 66.242043
$ws.Range("O6").Value = 0.1558818838066577
$ws.Range("P6").Value = 0.1558818838066577
$ws.Range("Q6").Value = 1.729498780233
$ws.Range("R6").Value = 15.565489022097
$ws.Range("S6").Value = 0.1403812174954666
$ws.Range("T6").Value = 0.1403812174954666
$ws.Range("I7").Value = 0.9005614640222286
$ws.Range("J7").Value = 0.9005614640222285
$ws.Range("M7").Value = 19.66149466666667
$ws.Range("N7").Value = 58.984484
$ws.Range("O7").Value = 0.1388032745500265
$ws.Range("P7").Value = 0.1388032745500265
$ws.Range("Q7").Value = 1.540012785092889
$ws.Range("R7").Value = 13.860115065836
$ws.Range("S7").Value = 0.1250008801398512
$ws.Range("T7").Value = 0.1250008801398512
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.008648666666666667
$ws.Range("H8").Value = 0.025946
$ws.Range("I8").Value = 0.0994385359777714
$ws.Range("J8").Value = 0.09943853597777139
$ws.Range("M8").Value = 36.81180933333333
$ws.Range("N8").Value = 110.435428
$ws.Range("O8").Value = 0.2598784967371026
$ws.Range("P8").Value = 0.2598784967371026
$ws.Range("Q8").Value = 0.3183730683208889
$ws.Range("R8").Value = 2.865357614888
$ws.Range("S8").Value = 0.02584193724764152
$ws.Range("T8").Value = 0.02584193724764152
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.008648666666666667
$ws.Range("H9").Value = 0.025946
$ws.Range("I9").Value = 0.0994385359777714
$ws.Range("J9").Value = 0.09943853597777139
$ws.Range("O9").Value = 0.1970278712683331
$ws.Range("P9").Value = 0.197027871268333
$ws.Range("Q9").Value = 0.2413757533155556
$ws.Range("R9").Value = 2.17238177984
$ws.Range("S9").Value = 0.01959216306573985
$ws.Range("T9").Value = 0.01959216306573984
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.008648666666666667
$ws.Range("H10").Value = 0.025946
$ws.Range("I10").Value = 0.0994385359777714
$ws.Range("J10").Value = 0.09943853597777139
$ws.Range("M10").Value = 21.95609833333333
$ws.Range("N10").Value = 65.868295
$ws.Range("O10").Value = 0.1550023737603119
$ws.Range("P10").Value = 0.1550023737603119
$ws.Range("Q10").Value = 0.1898909757855556
$ws.Range("R10").Value = 1.70901878207
$ws.Range("S10").Value = 0.01541320911980475
$ws.Range("T10").Value = 0.01541320911980474
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.008648666666666667
$ws.Range("H11").Value = 0.025946
$ws.Range("I11").Value = 0.0994385359777714
$ws.Range("J11").Value = 0.09943853597777139
$ws.Range("M11").Value = 13.23098133333333
$ws.Range("N11").Value = 39.692944
$ws.Range("O11").Value = 0.09340609987756826
$ws.Range("P11").Value = 0.09340609987756825
$ws.Range("Q11").Value = 0.1144303472248889
$ws.Range("R11").Value = 1.029873125024
$ws.Range("S11").Value = 0.009288165823218879
$ws.Range("T11").Value = 0.009288165823218878
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.008648666666666667
$ws.Range("H12").Value = 0.025946
$ws.Range("I12").Value = 0.0994385359777714
$ws.Range("J12").Value = 0.09943853597777139
$ws.Range("M12").Value = 22.080681
$ws.Range("N12").Value = 66.242043
$ws.Range("O12").Value = 0.1558818838066577
$ws.Range("P12").Value = 0.1558818838066577
$ws.Range("Q12").Value = 0.190968449742
$ws.Range("R12").Value = 1.718716047678
$ws.Range("S12").Value = 0.01550066631119111
$ws.Range("T12").Value = 0.01550066631119111
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.008648666666666667
$ws.Range("H13").Value = 0.025946
$ws.Range("I13").Value = 0.0994385359777714
$ws.Range("J13").Value = 0.09943853597777139
$ws.Range("M13").Value = 19.66149466666667
$ws.Range("N13").Value = 58.984484
$ws.Range("O13").Value = 0.1388032745500265
$ws.Range("P13").Value = 0.1388032745500265
$ws.Range("Q13").Value = 0.1700457135404445
$ws.Range("R13").Value = 1.530411421864
$ws.Range("S13").Value = 0.01380239441017529
$ws.Range("T13").Value = 0.01380239441017529
